$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-13 Sunday" "2025-07-14 Monday"

Replace-Text "223÷3=74, 1" "342÷6=57, 0"
Replace-Text "390÷9=43, 3" "759÷8=94, 7"
Replace-Text "481÷8=60, 1" "632÷5=126, 2"
Replace-Text "428÷6=71, 2" "937÷3=312, 1"
Replace-Text "222÷5=44, 2" "366÷4=91, 2"
Replace-Text "925÷7=132, 1" "223÷9=24, 7"
Replace-Text "278÷2=139, 0" "474÷7=67, 5"
Replace-Text "881÷3=293, 2" "663÷5=132, 3"
Replace-Text "909÷2=454, 1" "266÷3=88, 2"
Replace-Text "989÷5=197, 4" "863÷4=215, 3"
Replace-Text "615÷7=87, 6" "773÷6=128, 5"
Replace-Text "663÷9=73, 6" "209÷4=52, 1"
Replace-Text "533÷9=59, 2" "485÷6=80, 5"
Replace-Text "891÷8=111, 3" "485÷4=121, 1"
Replace-Text "577÷8=72, 1" "408÷6=68, 0"
Replace-Text "700÷8=87, 4" "169÷5=33, 4"
Replace-Text "523÷2=261, 1" "392÷5=78, 2"
Replace-Text "686÷2=343, 0" "935÷8=116, 7"
Replace-Text "168÷2=84, 0" "481÷4=120, 1"
Replace-Text "205÷5=41, 0" "240÷8=30, 0"
Replace-Text "242÷5=48, 2" "146÷4=36, 2"
Replace-Text "519÷7=74, 1" "894÷6=149, 0"
Replace-Text "980÷9=108, 8" "109÷6=18, 1"
Replace-Text "467÷3=155, 2" "228÷8=28, 4"
Replace-Text "441÷4=110, 1" "568÷4=142, 0"

Write-Host "Done"
